$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Row 2
$ws.Range("G2").Value = "16:12 - 2nd Half"

# Row 3
$ws.Range("G3").Value = "16:12 - 2nd Half"
$ws.Range("P3").Value = 23

# Row 7
$ws.Range("G7").Value = "16:12 - 2nd Half"
$ws.Range("P7").Value = 21

# Row 9
$ws.Range("G9").Value = "16:12 - 2nd Half"

# Row 12
$ws.Range("G12").Value = "16:12 - 2nd Half"

# Row 18
$ws.Range("G18").Value = "16:12 - 2nd Half"
$ws.Range("O18").Value = 3
$ws.Range("P18").Value = 21

# Row 21
$ws.Range("G21").Value = "16:12 - 2nd Half"
$ws.Range("H21").Value = 5
$ws.Range("J21").Value = 4
$ws.Range("K21").Value = 3
$ws.Range("P21").Value = 15

# Row 25
$ws.Range("G25").Value = "16:12 - 2nd Half"
$ws.Range("H25").Value = 14
$ws.Range("P25").Value = 22
$ws.Range("R25").Value = 6
$ws.Range("T25").Value = 3

# Row 28
$ws.Range("G28").Value = "16:12 - 2nd Half"
$ws.Range("H28").Value = 10
$ws.Range("I28").Value = 6
$ws.Range("P28").Value = 17
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = 5

# Row 36
$ws.Range("G36").Value = "16:12 - 2nd Half"
$ws.Range("J36").Value = 3
$ws.Range("P36").Value = 17
$ws.Range("R36").Value = 5

# Row 38
$ws.Range("D38").Value = "King Grace"
$ws.Range("E38").Value = "MSST"
$ws.Range("F38").Value = "MSST@SC"
$ws.Range("G38").Value = "Final"
$ws.Range("I38").Value = 6
$ws.Range("J38").Value = 2
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("P38").Value = 15
$ws.Range("R38").Value = 3
$ws.Range("T38").Value = 1
$ws.Range("U38").Value = 2
$ws.Range("V38").Value = 2

# Row 39
$ws.Range("D39").Value = "Chandler Bing"
$ws.Range("E39").Value = "VAN"
$ws.Range("F39").Value = "TENN@VAN"
$ws.Range("G39").Value = "16:12 - 2nd Half"
$ws.Range("H39").Value = 6
$ws.Range("J39").Value = 1
$ws.Range("N39").Value = 1
$ws.Range("O39").Value = 1
$ws.Range("P39").Value = 13
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 2
$ws.Range("T39").Value = 2
$ws.Range("U39").Value = 0
$ws.Range("V39").Value = 0

# Row 40
$ws.Range("D40").Value = "DeWayne Brown II"
$ws.Range("E40").Value = "TENN"
$ws.Range("G40").Value = "16:12 - 2nd Half"
$ws.Range("I40").Value = 4
$ws.Range("J40").Value = 3
$ws.Range("L40").Value = 1
$ws.Range("M40").Value = 1
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 2
$ws.Range("P40").Value = 18
$ws.Range("R40").Value = 5
$ws.Range("S40").Value = 0
$ws.Range("T40").Value = 0

# Row 45
$ws.Range("D45").Value = "AK Okereke"
$ws.Range("E45").Value = "VAN"
$ws.Range("F45").Value = "TENN@VAN"
$ws.Range("G45").Value = "16:12 - 2nd Half"
$ws.Range("I45").Value = 1
$ws.Range("K45").Value = 2
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("O45").Value = 3
$ws.Range("P45").Value = 10
$ws.Range("R45").Value = 0
$ws.Range("U45").Value = 1
$ws.Range("V45").Value = 2

# Row 46
$ws.Range("D46").Value = "Corey Chest"
$ws.Range("E46").Value = "MISS"
$ws.Range("F46").Value = "FLA@MISS"
$ws.Range("G46").Value = "Final"
$ws.Range("J46").Value = 1
$ws.Range("L46").Value = 1
$ws.Range("M46").Value = 2
$ws.Range("O46").Value = 1
$ws.Range("P46").Value = 12
$ws.Range("R46").Value = 2

# Row 47
$ws.Range("D47").Value = "Jayden Leverett"
$ws.Range("E47").Value = "VAN"
$ws.Range("F47").Value = "TENN@VAN"
$ws.Range("G47").Value = "16:12 - 2nd Half"
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 2
$ws.Range("P47").Value = 4
$ws.Range("R47").Value = 0
$ws.Range("T47").Value = 0
$ws.Range("U47").Value = 0
$ws.Range("V47").Value = 0

# Row 48
$ws.Range("D48").Value = "Max Smith"
$ws.Range("I48").Value = 2
$ws.Range("K48").Value = 1
$ws.Range("N48").Value = 0
$ws.Range("P48").Value = 5
$ws.Range("Q48").Value = 0
$ws.Range("T48").Value = 1
$ws.Range("V48").Value = 2

# Row 49
$ws.Range("D49").Value = "Niko Bundalo"
$ws.Range("I49").Value = 4
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 6
$ws.Range("R49").Value = 1
$ws.Range("T49").Value = 0
$ws.Range("U49").Value = 2
$ws.Range("V49").Value = 3

# Row 50
$ws.Range("D50").Value = "Zach Day"
$ws.Range("E50").Value = "MISS"
$ws.Range("F50").Value = "FLA@MISS"
$ws.Range("G50").Value = "Final"
$ws.Range("H50").Value = 3
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = 1
$ws.Range("L50").Value = 3
$ws.Range("N50").Value = 1
$ws.Range("O50").Value = 1
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = 4
$ws.Range("T50").Value = 3
$ws.Range("U50").Value = 0
$ws.Range("V50").Value = 0

# Row 54
$ws.Range("G54").Value = "16:12 - 2nd Half"

# Row 55
$ws.Range("G55").Value = "16:12 - 2nd Half"

# Row 60
$ws.Range("G60").Value = "16:12 - 2nd Half"

# OwnerTotals: The Oddities starter_pooh_total reflects Tyler Tanner pooh change (15 -> 14)
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Range("B4").Value = 14
